# "change to use testng.bat"
#
# - regression (sheet1): overwrite placeholder user/password values with
#   real ones (myUser/myPassword) and add a new scenario4 block.
# - stress (sheet2): add scenario3 rows for autoStart/RegRptFilePath and a
#   user/myPassword pair (carried over from the regression sheet's style).
# - logCompare (sheet3): drop the trailing scenario2 match/F43 rows.
# - status (sheet4): drop the trailing scenario2 rows.
# - three new sheets: batch, reconciliation, logon.
# - selections / active sheet bookkeeping to match the saved UI state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# regression
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("regression")

$ws1.Range("C9").Value = "myUser"
$ws1.Range("C10").Value = "myPassword"

$ws1.Range("A11").Value = "scenario4"
$ws1.Range("B11").Value = "batchFile"
$ws1.Range("C11").Value = "regress.bat"

$ws1.Range("A12").Value = "scenario4"
$ws1.Range("B12").Value = "user"
$ws1.Range("C12").Value = 'Peter;!@#$%''123'

$ws1.Range("A13").Value = "scenario4"
$ws1.Range("B13").Value = "password"
$ws1.Range("C13").Value = 'Tom;!@#$%''456'

# ---------------------------------------------------------------------
# stress
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("stress")

$ws2.Range("A6").Value = "scenario3"
$ws2.Range("B6").Value = "prntRegRpAtEnd"
$ws2.Range("C6").Value = "Y"

$ws2.Range("A7").Value = "scenario3"
$ws2.Range("B7").Value = "sprPort"
$ws2.Range("C7").Value = 1721

$ws2.Range("A8").Value = "scenario3"
$ws2.Range("B8").Value = "autoStart"
$ws2.Range("C8").Value = "C:\Users\abc.spr"

$ws2.Range("A9").Value = "scenario3"
$ws2.Range("B9").Value = "RegRptFilePath"
$ws2.Range("C9").Value = "C:\Work\RegReport.rtf"

$ws2.Range("A10").Value = "scenario3"
$ws2.Range("B10").Value = "user"
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("C10").Value = "myUser"
$ws2.Range("C10").NumberFormat = "@"

$ws2.Range("A11").Value = "scenario3"
$ws2.Range("B11").Value = "password"
$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("C11").Value = "myPassword"
$ws2.Range("C11").NumberFormat = "@"

# ---------------------------------------------------------------------
# logCompare - drop the trailing rows 12/13
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("logCompare")
$ws3.Range("A12:C13").ClearContents()

# ---------------------------------------------------------------------
# status - drop the trailing rows 3/4
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("status")
$ws4.Range("A3:C4").ClearContents()

# ---------------------------------------------------------------------
# batch (new sheet)
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $last)
$ws5.Name = "batch"
$ws5.Columns.Item(1).ColumnWidth = 17.59244791667
$ws5.Columns.Item(2).ColumnWidth = 16.30729166667
$ws5.Columns.Item(3).ColumnWidth = 20.73697916667

$ws5.Range("A1").Value = "Scenario No."
$ws5.Range("A1").NumberFormat = "@"
$ws5.Range("B1").Value = "key"
$ws5.Range("B1").NumberFormat = "@"
$ws5.Range("C1").Value = "value"
$ws5.Range("C1").NumberFormat = "@"

$ws5.Range("A2").Value = "scenario1"
$ws5.Range("B2").Value = "folderPath"
$ws5.Range("C2").Value = "C%3A%5CFINsim%5CINIfiles"

# ---------------------------------------------------------------------
# reconciliation (new sheet)
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6 = $wb.Worksheets.Add($null, $last)
$ws6.Name = "reconciliation"
$ws6.Columns.Item(1).ColumnWidth = 14.87760416667
$ws6.Columns.Item(2).ColumnWidth = 16.02213541667
$ws6.Columns.Item(3).ColumnWidth = 29.16666666667

$ws6.Range("A1").Value = "Scenario No."
$ws6.Range("A1").NumberFormat = "@"
$ws6.Range("B1").Value = "key"
$ws6.Range("B1").NumberFormat = "@"
$ws6.Range("C1").Value = "value"
$ws6.Range("C1").NumberFormat = "@"

$ws6.Range("A2").Value = "scenario1"
$ws6.Range("B2").Value = "sprAddress"
$ws6.Range("C2").Value = "127.0.0.1"

$ws6.Range("A3").Value = "scenario1"
$ws6.Range("B3").Value = "sprPort"
$ws6.Range("C3").Value = 1721

$ws6.Range("A4").Value = "scenario1"
$ws6.Range("B4").Value = "regressLogFile"
$ws6.Range("C4").Value = "C:\FINsim\regressLog.txt"

$ws6.Range("A5").Value = "scenario1"
$ws6.Range("B5").Value = "reconStart"
$ws6.Range("C5").Value = 1506081020

$ws6.Range("A6").Value = "scenario1"
$ws6.Range("B6").Value = "reconEnd"
$ws6.Range("C6").Value = 1506093045

$ws6.Range("A7").Value = "scenario1"
$ws6.Range("B7").Value = "reconLogFile"
$ws6.Range("C7").Value = "reconlog"

$ws6.Range("A8").Value = "scenario1"
$ws6.Range("B8").Value = "reconCsvFile"
$ws6.Range("C8").Value = "reconcsv"

$ws6.Range("A9").Value = "scenario1"
$ws6.Range("B9").Value = "reconrtfFile"
$ws6.Range("C9").Value = "recondoc"

$ws6.Range("A10").Value = "scenario1"
$ws6.Range("B10").Value = "reconDestFolder"
$ws6.Range("C10").Value = "C:\FINsim"

# ---------------------------------------------------------------------
# logon (new sheet)
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws7 = $wb.Worksheets.Add($null, $last)
$ws7.Name = "logon"
$ws7.Columns.Item(1).ColumnWidth = 12.30729166667
$ws7.Columns.Item(2).ColumnWidth = 16.16666666667
$ws7.Columns.Item(3).ColumnWidth = 17.73697916667

$ws7.Range("A1").Value = "Scenario No."
$ws7.Range("A1").NumberFormat = "@"
$ws7.Range("B1").Value = "key"
$ws7.Range("B1").NumberFormat = "@"
$ws7.Range("C1").Value = "value"
$ws7.Range("C1").NumberFormat = "@"

$ws7.Range("A2").Value = "scenario1"
$ws7.Range("B2").Value = "user"
$ws7.Range("B2").NumberFormat = "@"
$ws7.Range("C2").Value = "myUser"
$ws7.Range("C2").NumberFormat = "@"

$ws7.Range("A3").Value = "scenario1"
$ws7.Range("B3").Value = "password"
$ws7.Range("B3").NumberFormat = "@"
$ws7.Range("C3").Value = "myPassword"
$ws7.Range("C3").NumberFormat = "@"

$ws7.Range("A4").Value = "scenario1"
$ws7.Range("B4").Value = "sendLogonToAll"
$ws7.Range("C4").Value = "Y"
$ws7.Range("C4").NumberFormat = "@"

$ws7.Range("A5").Value = "scenario1"
$ws7.Range("B5").Value = "autoStart"
$ws7.Range("C5").Value = "C:\abc.spr"

# ---------------------------------------------------------------------
# selections - set per-sheet selection state; leave reconciliation for
# last so it ends up as the active/selected tab.
# ---------------------------------------------------------------------
$ws1.Range("J20").Select()
$ws2.Range("B10:C11").Select()
$ws3.Range("C17").Select()
$ws4.Range("I25").Select()
$ws5.Range("A2").Select()
$ws7.Range("C8").Select()
$ws6.Range("I24").Select()
